$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (student "Иванов Вячеслав") used to have a green-highlighted,
# partially filled-in row (C12:F12 = 2, G12/H12 blank). Bring it in line
# with the rest of the sheet (e.g. row 10): unfilled/bordered style,
# C12:H12 all = 5.
$ws.Range("C10:H10").Copy() | Out-Null
$ws.Range("C12:H12").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C12:H12").Value = 5

# K12 previously held the shared formula ("=SUM(C5:G5)" pattern inherited
# from K5); replace it with its own explicit SUM, matching K6/K10.
$ws.Range("K12").Formula = "=SUM(C12:H12)"

# Move the active selection from C21 to A8.
$ws.Range("A8").Select()
